# vmdefs.xlsx edit: "2 new machines AIA & 1 new test IA"
#
# Adds two new AHV-650 VMs (PRO_AIA, REC_AIA) to the first two data rows of the
# "vms" sheet (previously referencing the 651 cluster/datacenter), and fills in
# the formerly-blank 3rd data row with a new REC test "IA" machine. Also extends
# the params!sub_pe_lu650 subnet list (and its named range) by the 3 new
# AHV-650 network entries used by those rows.

$wb  = $excel.ActiveWorkbook
$vms = $wb.Worksheets.Item("vms")
$params = $wb.Worksheets.Item("params")

# ---------------------------------------------------------------------------
# params sheet: append the 3 new AHV-650 network names used by sub_pe_lu650
# (rows 25-41 already populated; 42/43 held leftover 651-network values in
# column D only -- now column B gets the matching 650-network values, and a
# brand new row 44 is added).
# ---------------------------------------------------------------------------
$params.Range("B42").Value = 'var.ahv_650_network["PRO_AIA_324"]'
$params.Range("B43").Value = 'var.ahv_650_network["REC_AIA_524"]'
$params.Range("B44").Value = 'var.ahv_650_network["REC_IA_547"]'

# Named range sub_pe_lu650 grows from B25:B41 to B25:B44 to cover the 3 new rows.
$wb.Names.Item("sub_pe_lu650").RefersTo = "=params!`$B`$25:`$B`$44"

# ---------------------------------------------------------------------------
# vms sheet, row 2: vsl-pro-aia-001 (was a DC3/651 dev IIS box, now a DC1/650
# PRO AIA box)
# ---------------------------------------------------------------------------
$vms.Range("B2").Value = "VSL-PRO-AIA-001"
$vms.Range("B2").Style = "Normal"
$vms.Range("C2").Value = "VSL-PRO-AIA-001"
$vms.Range("D2").Value = "VSL-PRO-AIA-001"
$vms.Range("E2").Value = "vsl-pro-aia-001"
$vms.Range("E2").Style = "Normal"
$vms.Range("F2").Value = "nutanix.dc1"
$vms.Range("G2").Value = "pe_lu650"
$vms.Range("H2").Value = "rhel8-dc1"
$vms.Range("J2").Value = 'var.ahv_650_network["PRO_AIA_324"]'
$vms.Range("L2").Value = 'var.ahv_650_storage["NUT_AHV_DC1_01"]'
$vms.Range("N2").Value = 24576
$vms.Range("O2").Value = 8
$vms.Range("R2").Value = "172.23.24.1 "
$vms.Range("R2").Style = "Normal"
$vms.Range("T2").Value = "172.23.24.253 "
$vms.Range("T2").Style = "Normal"

# ---------------------------------------------------------------------------
# vms sheet, row 3: vsl-rec-aia-001 (was a DC3/651 dev IDB box, now a DC1/650
# REC AIA box)
# ---------------------------------------------------------------------------
$vms.Range("B3").Value = "VSL-REC-AIA-001"
$vms.Range("C3").Value = "VSL-REC-AIA-001"
$vms.Range("D3").Value = "VSL-REC-AIA-001"
$vms.Range("E3").Value = "vsl-rec-aia-001"
$vms.Range("F3").Value = "nutanix.dc1"
$vms.Range("G3").Value = "pe_lu650"
$vms.Range("H3").Value = "rhel8-dc1"
$vms.Range("J3").Value = 'var.ahv_650_network["REC_AIA_524"]'
$vms.Range("L3").Value = 'var.ahv_650_storage["NUT_AHV_DC1_01"]'
$vms.Range("N3").Value = 24576
$vms.Range("O3").Value = 8
$vms.Range("R3").Value = "172.25.24.1"
$vms.Range("R3").Style = "Normal"
$vms.Range("T3").Value = "172.25.24.253 "
$vms.Range("T3").Style = "Normal"

# ---------------------------------------------------------------------------
# vms sheet, row 4: brand-new REC test "IA" machine (vsl-rec-iai-001), filling
# in what used to be a fully blank row.
# ---------------------------------------------------------------------------
$vms.Range("A4").Value = "LAN"
$vms.Range("B4").Value = "VSL-REC-IAI-001"
$vms.Range("C4").Value = "VSL-REC-IAI-001"
$vms.Range("D4").Value = "VSL-REC-IAI-001"
$vms.Range("E4").Value = "vsl-rec-iai-001"
$vms.Range("F4").Value = "nutanix.dc1"
$vms.Range("G4").Value = "pe_lu650"
$vms.Range("H4").Value = "rhel8-dc1"
$vms.Range("J4").Value = 'var.ahv_650_network["REC_IA_547"]'
$vms.Range("L4").Value = 'var.ahv_650_storage["NUT_AHV_DC1_01"]'
$vms.Range("N4").Value = 65535
$vms.Range("O4").Value = 24
$vms.Range("P4").Value = 1
$vms.Range("Q4").Value = 100
$vms.Range("R4").Value = "172.25.47.1 "
$vms.Range("R4").Style = "Normal"
$vms.Range("S4").Value = 24
$vms.Range("T4").Value = "172.25.47.253"
$vms.Range("T4").Style = "Normal"
$vms.Range("U4").Value = "DEV_TEST"

# ---------------------------------------------------------------------------
# Column E (HOSTNAME) on vms got wider to fit the new, longer hostnames.
# ---------------------------------------------------------------------------
$vms.Columns.Item(5).ColumnWidth = 24.140625

# ---------------------------------------------------------------------------
# View/selection cleanup: both sheets had scrolled/selected state left over
# from editing; reset to the plain view saved with the workbook.
# ---------------------------------------------------------------------------
$params.Activate()
$params.Range("B24").Select()

$vms.Activate()
$vms.Range("A1").Select()
